# Removing TWB PNPC table and adding two extra fields to TWB Episode instead.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- 1. Add the two new fields to "TWB Episodes" -------------------------
# Insert two blank columns at L:M. This pushes the existing
# "twb_previous_suicide_attempts" / "twb_method_of_suicide_attempt"
# columns (and their data) from L/M out to N/O automatically.
$twbEpisodes = $wb.Worksheets.Item("TWB Episodes")
$twbEpisodes.Range("L:M").Insert()

# Fill in the headers for the two newly-inserted columns.
$twbEpisodes.Range("L1").Value = "twb_primary_nominated_professional_contact_entry_date"
$twbEpisodes.Range("M1").Value = "twb_primary_nominated_professional_contact_exit_date"

# Fill in the data values for the two new columns.
$twbEpisodes.Range("L2").Value = 16042020
$twbEpisodes.Range("M2").Value = 9099999

$twbEpisodes.Range("L3").Value = 9099999
$twbEpisodes.Range("M3").Value = 9099999

# --- 2. Remove the now-unused "TWB PNPCs" table ---------------------------
$wb.Worksheets.Item("TWB PNPCs").Delete()
